# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-23 (serial 45192) to 2023-10-03 (serial 45202).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDateSerial = 45202

for ($row = 2; $row -le 306; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDateSerial
}
